# Excel DataTable scripts completed.
#
# 1) The EffectiveDate / PreviousExpDate columns (F and I) on customerInfo
#    held a stale "02172023" placeholder date for both data rows; bump it
#    to the corrected "02202023".
# 2) The workbook shipped two leftover helper sheets ("quoteChevron" and
#    "Staff") that are no longer needed now the DataTable scripts are
#    finished - remove them so only "customerInfo" remains.
# 3) Leave the sheet scrolled/selected the way the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customerInfo")

$ws.Range("F2").Value = "02202023"
$ws.Range("I2").Value = "02202023"
$ws.Range("F3").Value = "02202023"
$ws.Range("I3").Value = "02202023"

$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("quoteChevron").Delete()
[void]$wb.Worksheets.Item("Staff").Delete()
$excel.DisplayAlerts = $true

[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
[void]$ws.Range("H7").Select()
